$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("format")
$ws2 = $wb.Worksheets.Item("Tabelle1")

# Fix typo'd filename for NetCDF4-BEAM row
$ws1.Range("B3").Value = "netcdf4_beam.nc"

# Remove the "GeoTIFF-BigTIFF" row by shifting subsequent rows up,
# and append a new "HDF5" row at the end (binning/mosaic output formats)
$ws1.Range("A8").Value = "NetCDF4-CF"
$ws1.Range("B8").Value = "netcdf4_cf.nc"
$ws1.Range("A9").Value = "ENVI"
$ws1.Range("B9").Value = "envi"
$ws1.Range("A10").Value = "HDF5"
$ws1.Range("B10").Value = "hdf5.h5"

# Update the active selection to the last row, matching the post-edit state
$ws1.Activate()
$ws1.Range("A10:XFD10").Select()

# Sheet2 content is unaffected; just update the active selection
$ws2.Activate()
$ws2.Range("A12:XFD12").Select()

$ws1.Activate()
